$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.8606349999999999
$ws.Range("H2").Value = 2.581905
$ws.Range("I2").Value = 0.0262626340301864
$ws.Range("J2").Value = 0.0262626340301864
$ws.Range("M2").Value = 5.575746
$ws.Range("N2").Value = 16.727238
$ws.Range("O2").Value = 0.069238947264747
$ws.Range("P2").Value = 0.069238947264747
$ws.Range("Q2").Value = 4.798682158709999
$ws.Range("R2").Value = 43.18813942839
$ws.Range("S2").Value = 0.001818397132649426
$ws.Range("T2").Value = 0.001818397132649426
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.8606349999999999
$ws.Range("H3").Value = 2.581905
$ws.Range("I3").Value = 0.0262626340301864
$ws.Range("J3").Value = 0.0262626340301864
$ws.Range("O3").Value = 0.8150593598279631
$ws.Range("P3").Value = 0.815059359827963
$ws.Range("Q3").Value = 56.48859439385833
$ws.Range("R3").Value = 508.397349544725
$ws.Range("S3").Value = 0.0214056056800398
$ws.Range("T3").Value = 0.0214056056800398
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.8606349999999999
$ws.Range("H4").Value = 2.581905
$ws.Range("I4").Value = 0.0262626340301864
$ws.Range("J4").Value = 0.0262626340301864
$ws.Range("M4").Value = 9.317346333333333
$ws.Range("N4").Value = 27.952039
$ws.Range("O4").Value = 0.11570169290729
$ws.Range("P4").Value = 0.11570169290729
$ws.Range("Q4").Value = 8.018834361588333
$ws.Range("R4").Value = 72.16950925429499
$ws.Range("S4").Value = 0.00303863121749717
$ws.Range("T4").Value = 0.00303863121749717
$ws.Range("H5").Value = 58.40949000000001
$ws.Range("I5").Value = 0.5941299388474139
$ws.Range("J5").Value = 0.5941299388474139
$ws.Range("M5").Value = 5.575746
$ws.Range("N5").Value = 16.727238
$ws.Range("O5").Value = 0.069238947264747
$ws.Range("P5").Value = 0.069238947264747
$ws.Range("Q5").Value = 108.55882674318
$ws.Range("R5").Value = 977.0294406886201
$ws.Range("S5").Value = 0.04113693150426345
$ws.Range("T5").Value = 0.04113693150426345
$ws.Range("H6").Value = 58.40949000000001
$ws.Range("I6").Value = 0.5941299388474139
$ws.Range("J6").Value = 0.5941299388474139
$ws.Range("O6").Value = 0.8150593598279631
$ws.Range("P6").Value = 0.815059359827963
$ws.Range("S6").Value = 0.4842511676116
$ws.Range("T6").Value = 0.4842511676116
$ws.Range("H7").Value = 58.40949000000001
$ws.Range("I7").Value = 0.5941299388474139
$ws.Range("J7").Value = 0.5941299388474139
$ws.Range("M7").Value = 9.317346333333333
$ws.Range("N7").Value = 27.952039
$ws.Range("O7").Value = 0.11570169290729
$ws.Range("P7").Value = 0.11570169290729
$ws.Range("Q7").Value = 181.4071491611234
$ws.Range("R7").Value = 1632.66434245011
$ws.Range("S7").Value = 0.06874183973155047
$ws.Range("T7").Value = 0.06874183973155046
$ws.Range("G8").Value = 12.439858
$ws.Range("H8").Value = 37.319574
$ws.Range("I8").Value = 0.3796074271223998
$ws.Range("J8").Value = 0.3796074271223997
$ws.Range("M8").Value = 5.575746
$ws.Range("N8").Value = 16.727238
$ws.Range("O8").Value = 0.069238947264747
$ws.Range("P8").Value = 0.069238947264747
$ws.Range("Q8").Value = 69.36148848406799
$ws.Range("R8").Value = 624.2533963566119
$ws.Range("S8").Value = 0.02628361862783413
$ws.Range("T8").Value = 0.02628361862783412
$ws.Range("G9").Value = 12.439858
$ws.Range("H9").Value = 37.319574
$ws.Range("I9").Value = 0.3796074271223998
$ws.Range("J9").Value = 0.3796074271223997
$ws.Range("O9").Value = 0.8150593598279631
$ws.Range("P9").Value = 0.815059359827963
$ws.Range("Q9").Value = 816.5018769620033
$ws.Range("R9").Value = 7348.516892658029
$ws.Range("S9").Value = 0.3094025865363233
$ws.Range("T9").Value = 0.3094025865363232
$ws.Range("G10").Value = 12.439858
$ws.Range("H10").Value = 37.319574
$ws.Range("I10").Value = 0.3796074271223998
$ws.Range("J10").Value = 0.3796074271223997
$ws.Range("M10").Value = 9.317346333333333
$ws.Range("N10").Value = 27.952039
$ws.Range("O10").Value = 0.11570169290729
$ws.Range("P10").Value = 0.11570169290729
$ws.Range("Q10").Value = 115.9064653234873
$ws.Range("R10").Value = 1043.158187911386
$ws.Range("S10").Value = 0.04392122195824236
$ws.Range("T10").Value = 0.04392122195824235
